{"js": "// Document layout (in document order): a single title paragraph holding\n// the date, followed by a 20x5 table of arithmetic \"answer\" cells, each\n// cell containing its own one-paragraph body. `context.document.body\n// .paragraphs` walks the whole body in document order - title paragraph\n// first, then every table-cell paragraph row by row, column by column -\n// so a flat, position-indexed array lines up 1:1 with that collection\n// and lets us target each paragraph exactly (robust even though some\n// cell values, e.g. \"20+27=47\", repeat at more than one position).\nconst newValues = [\n  \"2023-06-05 Monday\",\n  \"33-8=25\",\n  \"40+44=84\",\n  \"59+4=63\",\n  \"30+46=76\",\n  \"38+12=50\",\n  \"28+8=36\",\n  \"74-2=72\",\n  \"24+22=46\",\n  \"81-36=45\",\n  \"71+15=86\",\n  \"64-2=62\",\n  \"4+28=32\",\n  \"19+55=74\",\n  \"28-20=8\",\n  \"98-10=88\",\n  \"95-35=60\",\n  \"71+16=87\",\n  \"13-13=0\",\n  \"68+20=88\",\n  \"43+35=78\",\n  \"6+91=97\",\n  \"68-47=21\",\n  \"63+18=81\",\n  \"41+46=87\",\n  \"81-42=39\",\n  \"20-18=2\",\n  \"70-13=57\",\n  \"78-46=32\",\n  \"34+24=58\",\n  \"63-35=28\",\n  \"99-10=89\",\n  \"19+15=34\",\n  \"18+41=59\",\n  \"95-37=58\",\n  \"47-27=20\",\n  \"54+0=54\",\n  \"34-30=4\",\n  \"95-59=36\",\n  \"52+10=62\",\n  \"88-62=26\",\n  \"2+80=82\",\n  \"1+79=80\",\n  \"94-41=53\",\n  \"43-2=41\",\n  \"79-37=42\",\n  \"62-55=7\",\n  \"63-34=29\",\n  \"65+1=66\",\n  \"47+29=76\",\n  \"21+68=89\",\n  \"30-30=0\",\n  \"35-30=5\",\n  \"93-51=42\",\n  \"28+57=85\",\n  \"19+40=59\",\n  \"83+10=93\",\n  \"63-14=49\",\n  \"77-71=6\",\n  \"61+13=74\",\n  \"12+82=94\",\n  \"95-39=56\",\n  \"2+41=43\",\n  \"83-51=32\",\n  \"53+32=85\",\n  \"11+9=20\",\n  \"47+17=64\",\n  \"3+43=46\",\n  \"27+67=94\",\n  \"18+12=30\",\n  \"68-52=16\",\n  \"42-41=1\",\n  \"39+30=69\",\n  \"17-9=8\",\n  \"76-8=68\",\n  \"66-62=4\",\n  \"73-49=24\",\n  \"73+23=96\",\n  \"54+38=92\",\n  \"89-37=52\",\n  \"14+9=23\",\n  \"58-41=17\",\n  \"33+28=61\",\n  \"80-60=20\",\n  \"77-58=19\",\n  \"75+1=76\",\n  \"8+10=18\",\n  \"62-16=46\",\n  \"81-72=9\",\n  \"68-15=53\",\n  \"49+37=86\",\n  \"27+10=37\",\n  \"85-35=50\",\n  \"25-3=22\",\n  \"26+34=60\",\n  \"28-7=21\",\n  \"10+69=79\",\n  \"63+32=95\",\n  \"92-69=23\",\n  \"8+14=22\",\n  \"89-31=58\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + paragraphs.items.length +\n    \" (expected \" + newValues.length + \")\"\n  );\n}\n\n// insertText(..., \"Replace\") rewrites just the paragraph's text run(s),\n// leaving the paragraph/run formatting (font, size, alignment) intact.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Document layout: a single title paragraph holding the date, followed by\n# a 20x5 table of arithmetic \"answer\" cells. We update the title text and\n# then every table cell by its (row, column) position - this is robust\n# even though some cell values (e.g. \"20+27=47\") repeat at more than one\n# position in the table.\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Item(1).Range.Text = '2023-06-05 Monday'\n\n# New value for each cell, laid out the same way as the table: 20 rows,\n# 5 columns each, in row-major order.\n$cellValues = @(\n    @('33-8=25', '40+44=84', '59+4=63', '30+46=76', '38+12=50'),\n    @('28+8=36', '74-2=72', '24+22=46', '81-36=45', '71+15=86'),\n    @('64-2=62', '4+28=32', '19+55=74', '28-20=8', '98-10=88'),\n    @('95-35=60', '71+16=87', '13-13=0', '68+20=88', '43+35=78'),\n    @('6+91=97', '68-47=21', '63+18=81', '41+46=87', '81-42=39'),\n    @('20-18=2', '70-13=57', '78-46=32', '34+24=58', '63-35=28'),\n    @('99-10=89', '19+15=34', '18+41=59', '95-37=58', '47-27=20'),\n    @('54+0=54', '34-30=4', '95-59=36', '52+10=62', '88-62=26'),\n    @('2+80=82', '1+79=80', '94-41=53', '43-2=41', '79-37=42'),\n    @('62-55=7', '63-34=29', '65+1=66', '47+29=76', '21+68=89'),\n    @('30-30=0', '35-30=5', '93-51=42', '28+57=85', '19+40=59'),\n    @('83+10=93', '63-14=49', '77-71=6', '61+13=74', '12+82=94'),\n    @('95-39=56', '2+41=43', '83-51=32', '53+32=85', '11+9=20'),\n    @('47+17=64', '3+43=46', '27+67=94', '18+12=30', '68-52=16'),\n    @('42-41=1', '39+30=69', '17-9=8', '76-8=68', '66-62=4'),\n    @('73-49=24', '73+23=96', '54+38=92', '89-37=52', '14+9=23'),\n    @('58-41=17', '33+28=61', '80-60=20', '77-58=19', '75+1=76'),\n    @('8+10=18', '62-16=46', '81-72=9', '68-15=53', '49+37=86'),\n    @('27+10=37', '85-35=50', '25-3=22', '26+34=60', '28-7=21'),\n    @('10+69=79', '63+32=95', '92-69=23', '8+14=22', '89-31=58'),\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $cellValues.Count; $r++) {\n    $row = $cellValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        # Setting Range.Text replaces just the cell's text content, so the\n        # existing run/paragraph formatting (font, size, alignment) is kept.\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
